$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.369.73'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '1.872.44'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '0.7140'
$ws.Range('E5').Value = '  +1.49%  '
$ws.Range('D6').Value = '241.34'
$ws.Range('E6').Value = '  +1.29%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '0.07892'
$ws.Range('E8').Value = '  -0.43%  '
$ws.Range('D9').Value = '0.3087'
$ws.Range('E9').Value = '  +1.48%  '
$ws.Range('D10').Value = '25.56'
$ws.Range('E10').Value = '  +4.12%  '
$ws.Range('D11').Value = '0.08250'
$ws.Range('E11').Value = '  +0.82%  '
$ws.Range('D12').Value = '0.7234'
$ws.Range('E12').Value = '  +0.64%  '
$ws.Range('D13').Value = '5.243'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').Value = '1.854.45'
$ws.Range('E14').Value = '  +9.73%  '
$ws.Range('D15').Value = '91.02'
$ws.Range('D16').Value = '29.363.46'
$ws.Range('E16').Value = '  +1.86%  '
$ws.Range('D17').Value = '5.850'
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('D18').Value = '244.42'
$ws.Range('E18').Value = '  +2.61%  '
$ws.Range('D19').Value = '0.000007824'
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').Value = '13.23'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').Value = '2.113.59'
$ws.Range('E21').Value = '  +10.49%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = '8.043'
$ws.Range('E23').Value = '  +6.63%  '
$ws.Range('D24').Value = '1.002'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = '0.1598'
$ws.Range('E25').Value = '  +11.90%  '
$ws.Range('D26').Value = '162.56'
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('D27').Value = '8.997'
$ws.Range('E27').Value = '  +1.17%  '
$ws.Range('D28').Value = '18.26'
$ws.Range('E28').Value = '  +0.91%  '
$ws.Range('E29').Value = '  -1.86%  '
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('D31').Value = '4.389'
$ws.Range('E31').Value = '  +1.39%  '
$ws.Range('D32').Value = '4.093'
$ws.Range('E33').Value = '  +0.20%  '
$ws.Range('D34').Value = '1.935'
$ws.Range('E34').Value = '  +0.89%  '
$ws.Range('E35').Value = '  +1.06%  '
$ws.Range('D36').Value = '0.7213'
$ws.Range('E36').Value = '  +1.06%  '
$ws.Range('D37').Value = '2.673'
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('D39').Value = '2.692'
$ws.Range('E39').Value = '  +0.20%  '
$ws.Range('D40').Value = '1.175.25'
$ws.Range('E40').Value = '  +1.35%  '
$ws.Range('D41').Value = '0.9063'
$ws.Range('E41').Value = '  -2.07%  '
$ws.Range('D42').Value = '6.105'
$ws.Range('E42').Value = '  +2.73%  '
$ws.Range('D43').Value = '72.53'
$ws.Range('E43').Value = '  +2.51%  '
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').Value = '102.14'
$ws.Range('E45').Value = '  +0.79%  '
$ws.Range('D46').Value = '0.5294'
$ws.Range('E46').Value = '  -0.40%  '
$ws.Range('D47').Value = '2.011.84'
$ws.Range('E47').Value = '  +10.27%  '
$ws.Range('D48').Value = '1.791'
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('D49').Value = '2.900'
$ws.Range('E49').Value = '  +5.79%  '
$ws.Range('D50').Value = '9.263'
$ws.Range('E50').Value = '  +0.78%  '
$ws.Range('D51').Value = '0.4286'
$ws.Range('E51').Value = '  +1.13%  '
